# Insert a new data row at row 23 (pushing the existing rows 23..104 down to 24..105),
# then populate the new row with the new Alcachofa price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value2 = 5
$ws.Range("B23").Value2 = "Macroferia Regional de Talca"
$ws.Range("C23").Value2 = "Maule"
$ws.Range("D23").Value2 = 44811
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value2 = 7
$ws.Range("F23").Value2 = 100112013
$ws.Range("G23").Value2 = "Alcachofa"
$ws.Range("H23").Value2 = "Madrigal"
$ws.Range("I23").Value2 = "Primera"
$ws.Range("J23").Value2 = 300
$ws.Range("K23").Value2 = 13000
$ws.Range("L23").Value2 = 13000
$ws.Range("M23").Value2 = 13000
$ws.Range("N23").Value2 = "`$/caja 40 unidades"
$ws.Range("O23").Value2 = "Provincia del Elquí"
$ws.Range("P23").Value2 = 325
$ws.Range("Q23").Value2 = 40
$ws.Range("R23").Value2 = "Hortaliza"
